# "Removed A4R slides from HPC & big data" (closes #121)
#
# The deck's 2nd slide (sldId 257) is a leftover, blank "Azure4Research"
# template slide with no shapes on it. Remove it so the deck goes straight
# from the title slide into the real content ("Big Data Analytics with
# HDInsight" / "Key learning objectives"). No other slide content changes.

$p = $ppt.ActivePresentation

$blankTemplateSlide = $p.Slides.FindBySlideID(257)
$blankTemplateSlide.Delete()
